$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Un-merge the old merged header cells (H1:L1 "Tackles", M1:P1 "Challenges", Q1:S1 "Blocks") ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- Row 1: replace the old two-tier / "Unnamed" header with the flattened header row ---
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Row 2 keeps the old header labels (unchanged) but is now hidden (kept for reference) ---
$ws.Rows.Item(2).Hidden = $true

# --- Row 3 is a blank separator row, also hidden ---
$ws.Rows.Item(3).Hidden = $true

# --- Clean the NaN "Tkl%" values (blank because TklW was 0) to an explicit 0 ---
$ws.Range("O5").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O19").Value = 0

# --- The totals row is hidden too ---
$ws.Rows.Item(20).Hidden = $true
